$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Materialize the previously-missing empty row 2 (no cell data, matches <row r="2"/>)
$ws.Rows("2:2").OutlineLevel = 0

# 2. Row 196: clear the (already-empty) trailing H196/I196 cells so they are dropped entirely
$ws.Range("H196:I196").ClearContents()

# 3. Append new incident rows 197-209
# Row 197
$ws.Range("A197").Value = "'2024-05-23"
$ws.Range("B197").Value = '09:49:58'
$ws.Range("C197").Value = '-'
$ws.Range("D197").Value = 'Cámara no detecta foam derecho'
$ws.Range("E197").Value = '-'
$ws.Range("F197").Value = '-'
$ws.Range("G197").Value = '-'
$ws.Range("H197").Value = '09:50:12'
$ws.Range("I197").Value = '0:00:14'

# Row 198
$ws.Range("A198").Value = "'2024-05-23"
$ws.Range("B198").Value = '09:50:39'
$ws.Range("C198").Value = '-'
$ws.Range("D198").Value = 'AOI (malla)'
$ws.Range("E198").Value = '-'
$ws.Range("F198").Value = '-'
$ws.Range("G198").Value = '-'
$ws.Range("H198").Value = '09:50:42'
$ws.Range("I198").Value = '0:00:03'

# Row 199
$ws.Range("A199").Value = "'2024-05-23"
$ws.Range("B199").Value = '10:39:08'
$ws.Range("C199").Value = '-'
$ws.Range("D199").Value = 'Cámara no detecta Pcb'
$ws.Range("E199").Value = '-'
$ws.Range("F199").Value = '-'
$ws.Range("G199").Value = '-'

# Row 200
$ws.Range("A200").Value = "'2024-05-23"
$ws.Range("B200").Value = '10:41:53'
$ws.Range("C200").Value = '-'
$ws.Range("D200").Value = 'Cámara no detecta Pcb'
$ws.Range("E200").Value = '-'
$ws.Range("F200").Value = '-'
$ws.Range("G200").Value = '-'
$ws.Range("H200").Value = '10:41:59'
$ws.Range("I200").Value = '0:00:06'

# Row 201
$ws.Range("A201").Value = "'2024-05-23"
$ws.Range("B201").Value = '10:42:21'
$ws.Range("C201").Value = '-'
$ws.Range("D201").Value = 'Cámara no detecta skeleton'
$ws.Range("E201").Value = '-'
$ws.Range("F201").Value = '-'
$ws.Range("G201").Value = '-'
$ws.Range("H201").Value = '10:42:30'
$ws.Range("I201").Value = '0:00:09'

# Row 202
$ws.Range("A202").Value = "'2024-05-23"
$ws.Range("B202").Value = '10:42:33'
$ws.Range("C202").Value = '-'
$ws.Range("D202").Value = 'No detecta presencia power CP'
$ws.Range("E202").Value = '-'
$ws.Range("F202").Value = '-'
$ws.Range("G202").Value = '-'
$ws.Range("H202").Value = '10:42:36'
$ws.Range("I202").Value = '0:00:03'

# Row 203
$ws.Range("A203").Value = "'2024-05-23"
$ws.Range("B203").Value = '10:46:15'
$ws.Range("C203").Value = '-'
$ws.Range("D203").Value = 'Cámara no detecta foam derecho'
$ws.Range("E203").Value = '-'
$ws.Range("F203").Value = '-'
$ws.Range("G203").Value = '-'
$ws.Range("H203").Value = '10:46:17'
$ws.Range("I203").Value = '0:00:02'

# Row 204
$ws.Range("A204").Value = "'2024-05-23"
$ws.Range("B204").Value = '10:48:22'
$ws.Range("C204").Value = '-'
$ws.Range("D204").Value = '-'
$ws.Range("E204").Value = '-'
$ws.Range("F204").Value = '-'
$ws.Range("G204").Value = 'Colisión placas'
$ws.Range("H204").Value = '10:48:25'
$ws.Range("I204").Value = '0:00:03'

# Row 205
$ws.Range("A205").Value = "'2024-05-23"
$ws.Range("B205").Value = '10:48:26'
$ws.Range("C205").Value = '-'
$ws.Range("D205").Value = '-'
$ws.Range("E205").Value = '-'
$ws.Range("F205").Value = '-'
$ws.Range("G205").Value = 'Soldadura defectuosa'
$ws.Range("H205").Value = '10:48:28'
$ws.Range("I205").Value = '0:00:02'

# Row 206
$ws.Range("A206").Value = "'2024-05-23"
$ws.Range("B206").Value = '11:36:30'
$ws.Range("C206").Value = '-'
$ws.Range("D206").Value = 'Cámara no detecta skeleton'
$ws.Range("E206").Value = '-'
$ws.Range("F206").Value = '-'
$ws.Range("G206").Value = '-'

# Row 207
$ws.Range("A207").Value = "'2024-05-23"
$ws.Range("B207").Value = '11:39:40'
$ws.Range("C207").Value = '-'
$ws.Range("D207").Value = 'Cámara no detecta skeleton'
$ws.Range("E207").Value = '-'
$ws.Range("F207").Value = '-'
$ws.Range("G207").Value = '-'
$ws.Range("H207").Value = '11:40:25'
$ws.Range("I207").Value = '0:00:45'

# Row 208
$ws.Range("A208").Value = "'2024-05-23"
$ws.Range("B208").Value = '11:57:13'
$ws.Range("C208").Value = 'Etiquetadora'
$ws.Range("D208").Value = '-'
$ws.Range("E208").Value = '-'
$ws.Range("F208").Value = '-'
$ws.Range("G208").Value = '-'
$ws.Range("H208").Value = '11:57:24'
$ws.Range("I208").Value = '0:00:11'

# Row 209
$ws.Range("A209").Value = "'2024-05-23"
$ws.Range("B209").Value = '11:57:20'
$ws.Range("C209").Value = 'Etiquetadora'
$ws.Range("D209").Value = '-'
$ws.Range("E209").Value = '-'
$ws.Range("F209").Value = '-'
$ws.Range("G209").Value = '-'
$ws.Range("H209").Value = '11:57:24'
$ws.Range("I209").Value = '0:00:04'

# The apostrophe-prefix trick above marks A197:A209 with a "quote prefix"
# cell style; reset back to the Normal style so the cells keep plain General
# formatting (matching the rest of the sheet) while the stored value stays text.
$ws.Range("A197:A209").Style = "Normal"

